$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Select the DP table range on the original sheet before duplicating it
$ws1.Range("B3:L13").Select()

# Duplicate the sheet (Move or Copy -> Create a copy), inserting after itself
$ws1.Copy($null, $ws1)

$ws2 = $wb.Worksheets.Item(2)

# On the new copy, delete column C ("Beneficio") entirely, shifting cells left
$ws2.Range("C1").EntireColumn.Delete()

# Clear the helper row (former row 2 "j" header) and the DP table body
$ws2.Range("B2:K2").ClearContents()
$ws2.Range("C5:K13").ClearContents()

$ws2.Range("M16").Select()
